# Updated cryptos list (price/volume refresh), matching the GitHub Actions commit.
# D-column "Price" values are forced to text (leading apostrophe) so Excel
# doesn't silently reinterpret numeric-looking strings (e.g. "211.30",
# "1.00") as floating point numbers and strip the formatting/trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.714.57"
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = "'1.600.92"
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = "'211.30"
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('E6').Value = '  -0.45%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').Value = "'19.67"
$ws.Range('E10').Value = '  +0.93%  '
$ws.Range('D11').Value = "'0.0845"
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('D12').Value = "'1.824.96"
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').Value = "'1.602.84"
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').Value = "'26.684.49"
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = "'0.0₃0741"
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('D20').Value = "'7.21"
$ws.Range('E20').Value = '  +2.87%  '
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('D22').Value = "'4.29"
$ws.Range('E22').Value = '  +0.58%  '
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').Value = "'8.96"
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('D25').Value = "'144.17"
$ws.Range('E25').Value = '  -0.77%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('E28').Value = '  -0.88%  '
$ws.Range('D29').Value = "'15.35"
$ws.Range('E29').Value = '  +0.63%  '
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('E32').Value = '  +1.46%  '
$ws.Range('D33').Value = "'2.98"
$ws.Range('E33').Value = '  +1.24%  '
$ws.Range('D34').Value = "'1.295.56"
$ws.Range('E34').Value = '  +1.40%  '
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('E36').Value = '  +0.95%  '
$ws.Range('D37').Value = "'0.603"
$ws.Range('E37').Value = '  -2.10%  '
$ws.Range('E38').Value = '  +15.13%  '
$ws.Range('E39').Value = '  -0.63%  '
$ws.Range('D40').Value = "'0.823"
$ws.Range('E40').Value = '  -1.88%  '
$ws.Range('E41').Value = '  -1.39%  '
$ws.Range('D42').Value = "'2.20"
$ws.Range('E42').Value = '  -0.31%  '
$ws.Range('D43').Value = "'0.779"
$ws.Range('E43').Value = '  -0.61%  '
$ws.Range('D44').Value = "'63.10"
$ws.Range('E44').Value = '  -1.68%  '
$ws.Range('D45').Value = "'1.738.06"
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('E46').Value = '  +0.55%  '
$ws.Range('E47').Value = '  -2.90%  '
$ws.Range('E48').Value = '  -0.71%  '
$ws.Range('D49').Value = "'0.0518"
$ws.Range('E49').Value = '  +2.02%  '
$ws.Range('E50').Value = '  -0.07%  '
$ws.Range('E51').Value = '  -0.51%  '
